# Apply the "T1..T14 attendance sheet" edit to the Ninja worksheet and
# update the active sheet / selection state on both the Astronauta and
# Ninja worksheets, mirroring the author's manual Excel session.

$wb = $excel.ActiveWorkbook

$wsAstro = $wb.Worksheets.Item("Astronauta")
$wsNinja = $wb.Worksheets.Item("Ninja")

# --- Fill in column B (the "T1" attendance/participation column) for
#     rows 2..21 on the Ninja sheet. This drives the existing P column
#     IFERROR(SUM(..)/COUNT(..)*100,0) formulas from 0 to 100 wherever a
#     1 was entered.
$attendance = @(1, 1, 1, 0, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $attendance.Length; $i++) {
    $row = $i + 2
    $wsNinja.Cells.Item($row, 2).Value = $attendance[$i]
}

# --- Selection / active-sheet bookkeeping, matching the diff:
#     Astronauta loses tabSelected and its selection moves to B2.
$wsAstro.Activate()
$wsAstro.Range("B2").Select()

#     Ninja becomes the tab-selected / active sheet, with selection B6.
$wsNinja.Activate()
$wsNinja.Range("B6").Select()
